$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = [double]"964.88884562729652"
$ws.Range("H2").Value = [double]"74.811264805230593"
$ws.Range("I2").Value = [double]"13.145267620013534"
$ws.Range("J2").Value = [double]"1681.8927482316692"
$ws.Range("K2").Value = [double]"-3.8316688590874983E-10"
$ws.Range("O2").Value = [double]"13.145267620013534"

$ws.Range("G3").Value = [double]"678.73515526746507"
$ws.Range("H3").Value = [double]"81.730322542628571"
$ws.Range("I3").Value = [double]"1.0928213395648341"
$ws.Range("J3").Value = [double]"970.59278274442465"
$ws.Range("K3").Value = [double]"1.0419807390234324E-8"
$ws.Range("O3").Value = [double]"1.0928213395648341"

$ws.Range("G4").Value = [double]"915.7355279937052"
$ws.Range("H4").Value = [double]"99.61347269596277"
$ws.Range("I4").Value = [double]"3.1767368924536519"
$ws.Range("J4").Value = [double]"970.59279007387033"
$ws.Range("K4").Value = [double]"1.0419807505105544E-8"
$ws.Range("O4").Value = [double]"3.1767368924536519"

$ws.Range("G5").Value = [double]"568.82366691035327"
$ws.Range("H5").Value = [double]"50.850852959614393"
$ws.Range("I5").Value = [double]"3.5713104474987101"
$ws.Range("J5").Value = [double]"972.254092343146"
$ws.Range("K5").Value = [double]"-5.9805267946958673E-10"
$ws.Range("O5").Value = [double]"3.5713104474987101"

$ws.Range("G6").Value = [double]"964.88882169024339"
$ws.Range("H6").Value = [double]"74.811888682028169"
$ws.Range("I6").Value = [double]"13.400003353681761"
$ws.Range("J6").Value = [double]"1273.7205982263815"
$ws.Range("K6").Value = [double]"-3.2827068302375881E-10"
$ws.Range("O6").Value = [double]"13.400003353681761"

$ws.Range("G7").Value = [double]"421.76128744832067"
$ws.Range("H7").Value = [double]"53.834249902219547"
$ws.Range("I7").Value = [double]"2.1681142742195001"
$ws.Range("J7").Value = [double]"970.59278563297494"
$ws.Range("K7").Value = [double]"1.0419807339439452E-8"
$ws.Range("O7").Value = [double]"2.1681142742195001"

$ws.Range("G8").Value = [double]"278.49822985129339"
$ws.Range("H8").Value = [double]"26.297106843548995"
$ws.Range("I8").Value = [double]"2.3347051051593266"
$ws.Range("J8").Value = [double]"970.59278625091144"
$ws.Range("K8").Value = [double]"1.0419807329654527E-8"
$ws.Range("O8").Value = [double]"2.3347051051593266"

$ws.Range("G9").Value = [double]"964.88869653406766"
$ws.Range("H9").Value = [double]"74.813682576139243"
$ws.Range("I9").Value = [double]"11.749063647516842"
$ws.Range("J9").Value = [double]"1555.5651822845109"
$ws.Range("K9").Value = [double]"2.0197099284881638E-10"
$ws.Range("O9").Value = [double]"11.749063647516842"

$ws.Range("G11").Value = [double]"964.88854770721957"
$ws.Range("H11").Value = [double]"74.815033377504776"
$ws.Range("I11").Value = [double]"3.0495049274861405"
$ws.Range("J11").Value = [double]"1057.5532878325171"
$ws.Range("K11").Value = [double]"-1.4166384473778876E-9"
$ws.Range("O11").Value = [double]"3.0495049274861405"

$ws.Range("G12").Value = [double]"473.28884901725036"
$ws.Range("H12").Value = [double]"45.092370491185598"
$ws.Range("I12").Value = [double]"0.39249235889455036"
$ws.Range("J12").Value = [double]"970.59278188751932"
$ws.Range("K12").Value = [double]"1.0419807396522696E-8"
$ws.Range("O12").Value = [double]"0.39249235889455036"

$ws.Range("G13").Value = [double]"964.88871178035379"
$ws.Range("H13").Value = [double]"74.812952331288031"
$ws.Range("I13").Value = [double]"9.9368399301019323"
$ws.Range("J13").Value = [double]"1124.7115291208486"
$ws.Range("K13").Value = [double]"-1.4232604389657126E-9"
$ws.Range("O13").Value = [double]"9.9368399301019323"

$ws.Range("G14").Value = [double]"964.88928064066693"
$ws.Range("H14").Value = [double]"74.805661459739568"
$ws.Range("I14").Value = [double]"20.181941844009138"
$ws.Range("J14").Value = [double]"2764.9691170165443"
$ws.Range("K14").Value = [double]"9.48539401295092E-11"
$ws.Range("O14").Value = [double]"20.181941844009138"

$ws.Range("G15").Value = [double]"964.88864767771622"
$ws.Range("H15").Value = [double]"74.813759832252344"
$ws.Range("I15").Value = [double]"7.9525041612801468"
$ws.Range("J15").Value = [double]"1338.9802235312629"
$ws.Range("K15").Value = [double]"-9.1248638888952726E-11"
$ws.Range("O15").Value = [double]"7.9525041612801468"

$ws.Range("G16").Value = [double]"421.76129154565768"
$ws.Range("H16").Value = [double]"53.834247425615516"
$ws.Range("I16").Value = [double]"1.6438156844787695"
$ws.Range("J16").Value = [double]"970.59278398657909"
$ws.Range("K16").Value = [double]"1.0419807430241384E-8"
$ws.Range("O16").Value = [double]"1.6438156844787695"

$ws.Range("G17").Value = [double]"616.0446778605924"
$ws.Range("H17").Value = [double]"64.774596627603202"
$ws.Range("I17").Value = [double]"2.4736579087627897"
$ws.Range("J17").Value = [double]"970.59278680135913"
$ws.Range("K17").Value = [double]"1.0419807407500243E-8"
$ws.Range("O17").Value = [double]"2.4736579087627897"

$ws.Range("G18").Value = [double]"964.88880505448697"
$ws.Range("H18").Value = [double]"74.811724470121632"
$ws.Range("I18").Value = [double]"12.148931852528131"
$ws.Range("J18").Value = [double]"1362.2980919322511"
$ws.Range("K18").Value = [double]"2.2401638874640371E-10"
$ws.Range("O18").Value = [double]"12.148931852528131"

$ws.Range("G19").Value = [double]"964.88918494834752"
$ws.Range("H19").Value = [double]"74.807408220819568"
$ws.Range("I19").Value = [double]"19.622175803443472"
$ws.Range("J19").Value = [double]"2763.1200801819646"
$ws.Range("K19").Value = [double]"9.2573334085842576E-11"
$ws.Range("O19").Value = [double]"19.622175803443472"

$ws.Range("G20").Value = [double]"964.88876465564817"
$ws.Range("H20").Value = [double]"74.812690766102605"
$ws.Range("I20").Value = [double]"12.473148749000446"
$ws.Range("J20").Value = [double]"2713.9457906656512"
$ws.Range("K20").Value = [double]"1.0162290173411643E-10"
$ws.Range("O20").Value = [double]"12.473148749000446"

$ws.Range("G21").Value = [double]"964.88864323205053"
$ws.Range("H21").Value = [double]"74.813795258566671"
$ws.Range("I21").Value = [double]"7.7242545641359159"
$ws.Range("J21").Value = [double]"1436.1224525572404"
$ws.Range("K21").Value = [double]"-1.1671319686427043E-9"
$ws.Range("O21").Value = [double]"7.7242545641359159"

$ws.Range("G22").Value = [double]"964.88854965561848"
$ws.Range("H22").Value = [double]"74.815009967317621"
$ws.Range("I22").Value = [double]"3.2194998734757645"
$ws.Range("J22").Value = [double]"1161.4308784227276"
$ws.Range("K22").Value = [double]"5.484154367584804E-10"
$ws.Range("O22").Value = [double]"3.2194998734757645"

$ws.Range("G23").Value = [double]"964.88857648858027"
$ws.Range("H23").Value = [double]"74.814713927879978"
$ws.Range("I23").Value = [double]"5.2815531048241402"
$ws.Range("J23").Value = [double]"1104.8488931079476"
$ws.Range("K23").Value = [double]"9.657426235165033E-10"
$ws.Range("O23").Value = [double]"5.2815531048241402"

$ws.Range("G24").Value = [double]"655.0980073890463"
$ws.Range("H24").Value = [double]"90.271617499368674"
$ws.Range("I24").Value = [double]"4.1609241773326486"
$ws.Range("J24").Value = [double]"1084.345043995366"
$ws.Range("K24").Value = [double]"7.0484870460466807E-10"
$ws.Range("O24").Value = [double]"4.1609241773326486"

$ws.Range("G25").Value = [double]"964.8887965567028"
$ws.Range("H25").Value = [double]"74.811852171993678"
$ws.Range("I25").Value = [double]"11.998777589526354"
$ws.Range("J25").Value = [double]"1954.0513283830689"
$ws.Range("K25").Value = [double]"-8.1452846636326227E-11"
$ws.Range("O25").Value = [double]"11.998777589526354"

$ws.Range("G26").Value = [double]"933.99325073901957"
$ws.Range("H26").Value = [double]"77.491038248343827"
$ws.Range("I26").Value = [double]"5.8916905519731104"
$ws.Range("J26").Value = [double]"1015.1743177040878"
$ws.Range("K26").Value = [double]"-5.4153244412746727E-10"
$ws.Range("O26").Value = [double]"5.8916905519731104"

$ws.Range("G27").Value = [double]"964.88854071567107"
$ws.Range("H27").Value = [double]"74.81512460443102"
$ws.Range("I27").Value = [double]"3.0896506469680207"
$ws.Range("J27").Value = [double]"970.59278962432313"
$ws.Range("K27").Value = [double]"1.0419807424657743E-8"
$ws.Range("O27").Value = [double]"3.0896506469680207"

$ws.Range("G28").Value = [double]"964.88854827156547"
$ws.Range("H28").Value = [double]"74.815032559851218"
$ws.Range("I28").Value = [double]"3.2034031541909664"
$ws.Range("J28").Value = [double]"970.59279021413749"
$ws.Range("K28").Value = [double]"1.0419807374486287E-8"
$ws.Range("O28").Value = [double]"3.2034031541909664"

$ws.Range("G29").Value = [double]"964.88872728460217"
$ws.Range("H29").Value = [double]"74.81296592386127"
$ws.Range("I29").Value = [double]"10.970419906812999"
$ws.Range("J29").Value = [double]"1944.924695456825"
$ws.Range("K29").Value = [double]"-1.4778067407646721E-10"
$ws.Range("O29").Value = [double]"10.970419906812999"

$ws.Range("G30").Value = [double]"186.87270558583845"
$ws.Range("H30").Value = [double]"40.180913491260469"
$ws.Range("I30").Value = [double]"1.5182151582236156"
$ws.Range("J30").Value = [double]"970.59278365941088"
$ws.Range("K30").Value = [double]"1.0419807436529756E-8"
$ws.Range("O30").Value = [double]"1.5182151582236156"

$ws.Range("G31").Value = [double]"933.99325046175659"
$ws.Range("H31").Value = [double]"77.49103871671052"
$ws.Range("I31").Value = [double]"4.7707857025662772"
$ws.Range("J31").Value = [double]"997.10009259745311"
$ws.Range("K31").Value = [double]"-3.6702085090052558E-10"
$ws.Range("O31").Value = [double]"4.7707857025662772"

$ws.Range("G32").Value = [double]"278.49822395790767"
$ws.Range("H32").Value = [double]"26.297117588981369"
$ws.Range("I32").Value = [double]"1.5023040163036467"
$ws.Range("J32").Value = [double]"970.59278361985639"
$ws.Range("K32").Value = [double]"1.0419807463228235E-8"
$ws.Range("O32").Value = [double]"1.5023040163036467"

$ws.Range("G33").Value = [double]"964.88855880233905"
$ws.Range("H33").Value = [double]"74.814919458513899"
$ws.Range("I33").Value = [double]"4.1614056656251428"
$ws.Range("J33").Value = [double]"1382.6695706314126"
$ws.Range("K33").Value = [double]"1.2465416373625834E-9"
$ws.Range("O33").Value = [double]"4.1614056656251428"

$ws.Range("G34").Value = [double]"964.88864074037622"
$ws.Range("H34").Value = [double]"74.814149158790897"
$ws.Range("I34").Value = [double]"9.1883255227452736"
$ws.Range("J34").Value = [double]"1251.9442881198179"
$ws.Range("K34").Value = [double]"4.70070796307977E-10"
$ws.Range("O34").Value = [double]"9.1883255227452736"

$ws.Range("G35").Value = [double]"757.74013105488939"
$ws.Range("H35").Value = [double]"86.869471052109333"
$ws.Range("I35").Value = [double]"3.715525875401787"
$ws.Range("J35").Value = [double]"1074.5932768852201"
$ws.Range("K35").Value = [double]"1.3871111094382881E-10"
$ws.Range("O35").Value = [double]"3.715525875401787"

$ws.Range("G36").Value = [double]"964.88907219871851"
$ws.Range("H36").Value = [double]"74.808680221694246"
$ws.Range("I36").Value = [double]"17.711728456662104"
$ws.Range("J36").Value = [double]"1607.7397613682006"
$ws.Range("K36").Value = [double]"-1.0589158246699423E-9"
$ws.Range("O36").Value = [double]"17.711728456662104"

$ws.Range("G37").Value = [double]"616.04468087569649"
$ws.Range("H37").Value = [double]"64.774595349560698"
$ws.Range("I37").Value = [double]"2.685779143800302"
$ws.Range("J37").Value = [double]"1046.6031240474372"
$ws.Range("K37").Value = [double]"1.2651717421156755E-9"
$ws.Range("O37").Value = [double]"2.685779143800302"

$ws.Range("G38").Value = [double]"964.88897038019832"
$ws.Range("H38").Value = [double]"74.809928497331711"
$ws.Range("I38").Value = [double]"15.981823243407925"
$ws.Range("J38").Value = [double]"1680.9620413392272"
$ws.Range("K38").Value = [double]"-4.0455244104217884E-11"
$ws.Range("O38").Value = [double]"15.981823243407925"

$ws.Range("G39").Value = [double]"964.88857446814643"
$ws.Range("H39").Value = [double]"74.814844591335628"
$ws.Range("I39").Value = [double]"6.460074738597382"
$ws.Range("J39").Value = [double]"1023.0688600437626"
$ws.Range("K39").Value = [double]"7.427039826830365E-10"
$ws.Range("O39").Value = [double]"6.460074738597382"
